$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.142.82'
$ws.Range('E2').Value = '  -4.29%  '
$ws.Range('D3').Value = '3.311.67'
$ws.Range('E3').Value = '  -5.91%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '560.08'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.90%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '181.67'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -5.66%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  -2.26%  '
$ws.Range('D9').Value = '3.303.05'
$ws.Range('E9').Value = '  -5.81%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.190'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -6.51%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.590'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -4.38%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '47.85'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -8.00%  '
$ws.Range('E13').Value = '  -6.45%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '639.87'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.50%  '
$ws.Range('E15').Value = '  -6.21%  '
$ws.Range('D16').Value = '3.841.28'
$ws.Range('E16').Value = '  -5.79%  '
$ws.Range('D17').Value = '66.148.63'
$ws.Range('E17').Value = '  -4.31%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '17.96'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.59%  '
$ws.Range('E19').Value = '  -3.57%  '
$ws.Range('D20').Value = '3.305.19'
$ws.Range('E20').Value = '  -6.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.50'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -7.63%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.907'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.71%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '17.64'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.68%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '107.08'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +5.36%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.06'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -7.12%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.03'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -7.18%  '
$ws.Range('E27').Value = '  -0.47%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.70'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -6.73%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.58'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.87%  '
$ws.Range('E30').Value = '  -6.98%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '30.80'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -6.01%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.06'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.59%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.39'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.80%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.11'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.52%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '549.59'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +8.98%  '
$ws.Range('E36').Value = '  -3.73%  '
$ws.Range('B37').Value = 'Dai'
$ws.Range('C37').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.999'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.13%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '57.41'
$ws.Range('D38').Style = 'Normal'
$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').Value = '3.724.55'
$ws.Range('E39').Value = '  +0.24%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.53'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.47%  '
$ws.Range('B41').Value = 'PEPE'
$ws.Range('C41').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D41').Value = '0.0₃0719'
$ws.Range('E41').Value = '  -9.52%  '
$ws.Range('B42').Value = 'Fetch.AI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.74'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -6.19%  '
$ws.Range('E43').Value = '  -3.90%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.37'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +24.95%  '
$ws.Range('B45').Value = 'TheGraph'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.343'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -6.54%  '
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '32.31'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -6.36%  '
$ws.Range('E47').Value = '  -5.67%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.25'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.55%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.64'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -6.69%  '
$ws.Range('E50').Value = '  -3.56%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.998'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.22%  '
